# Revert "Adding the RES Hourly Production Forecast to the Portfolio"
#
# 1) Roll the Timestamp column (A2:A97, one day's worth of 15-minute
#    readings) back from 2024-09-24 to 2024-08-29.
# 2) Restore the pre-revert forecast figures (columns B:F) for rows 29-58.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Timestamps: 2024-09-24 -> 2024-08-29, same time-of-day, row 2..97 ---
$startDate = Get-Date -Year 2024 -Month 8 -Day 29 -Hour 0 -Minute 0 -Second 0
for ($i = 0; $i -lt 96; $i++) {
    $row = 2 + $i
    $ts = $startDate.AddMinutes(15 * $i)
    $ws.Range("A$row").Value = $ts.ToString("yyyy-MM-dd HH:mm:ss") + "+00:00"
}

# --- 2) Forecast values: columns B (Power_MW), C (Next_Power_MW),
#        D (Average_Power_MW... lag), E, F (Energy_MWh) for rows 29-58 ---
$forecastRows = @(
    @{ Row=29; B='1765.557335535685'; C='0.0017655573355356'; E='0.0008827786677678'; F='0.0002206946669419' },
    @{ Row=30; B='28658.85710652669'; C='0.0286588571065266'; D='0.0017655573355356'; E='0.0152122072210311'; F='0.0038030518052577' },
    @{ Row=31; B='56037.27966308594'; C='0.0560372796630859'; D='0.0286588571065266'; E='0.0423480683848063'; F='0.0105870170962015' },
    @{ Row=32; B='48097.56840006511'; C='0.0480975684000651'; D='0.0560372796630859'; E='0.0520674240315754'; F='0.0130168560078938' },
    @{ Row=33; B='83071.54130045572'; C='0.0830715413004556'; D='0.0480975684000651'; E='0.0655845548502603'; F='0.016396138712565' },
    @{ Row=34; B='119456.5695800781'; C='0.1194565695800781'; D='0.0830715413004556'; E='0.1012640554402668'; F='0.0253160138600666' },
    @{ Row=35; B='194109.4351399739'; C='0.1941094351399739'; D='0.1194565695800781'; E='0.1567830023600259'; F='0.0391957505900064' },
    @{ Row=36; B='603478.037923177'; C='0.6034780379231771'; D='0.1941094351399739'; E='0.3987937365315755'; F='0.09969843413289381' },
    @{ Row=37; B='406300.1061197916'; C='0.4063001061197916'; D='0.6034780379231771'; E='0.5048890720214844'; F='0.126222268005371' },
    @{ Row=38; B='486426.8177083333'; C='0.4864268177083333'; D='0.4063001061197916'; E='0.4463634619140624'; F='0.1115908654785156' },
    @{ Row=39; B='588672.568359375'; C='0.588672568359375'; D='0.4864268177083333'; E='0.5375496930338541'; F='0.1343874232584635' },
    @{ Row=40; B='658117.0885416666'; C='0.6581170885416666'; D='0.588672568359375'; E='0.6233948284505209'; F='0.1558487071126302' },
    @{ Row=41; B='657008.0100911459'; C='0.6570080100911458'; D='0.6581170885416666'; E='0.6575625493164062'; F='0.1643906373291015' },
    @{ Row=42; B='532883.3001302084'; C='0.5328833001302083'; D='0.6570080100911458'; E='0.5949456551106771'; F='0.1487364137776692' },
    @{ Row=43; B='0'; C='0'; D='0.5328833001302083'; E='0.2664416500651042'; F='0.066610412516276' },
    @{ Row=44; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=45; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=46; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=47; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=48; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=49; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=50; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=51; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=52; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=53; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=54; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=55; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=56; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=57; B='0'; C='0'; D='0'; E='0'; F='0' },
    @{ Row=58; D='0'; E='0'; F='0' }
)

foreach ($r in $forecastRows) {
    $row = $r.Row
    foreach ($col in @('B', 'C', 'D', 'E', 'F')) {
        if ($r.ContainsKey($col)) {
            $ws.Range("$col$row").Value = [double]$r[$col]
        }
    }
}
